$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-211 correspond to Generation 0..209 -> Fitness 7310
$ws.Range("C2:C211").Value = 7310

# Rows 212-252 correspond to Generation 210..250 -> Fitness 7293
$ws.Range("C212:C252").Value = 7293
